$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timesheet entries for row 11 (minutes/hours split) and row 12/13 minutes
$ws.Range("D11").Value = 1
$ws.Range("F11").Value = 15
$ws.Range("F12").Value = 30
$ws.Range("F13").Value = 45

# Update the active selection to reflect where the user left off editing
$ws.Range("D13").Select()
